# The deck's single custom "Integral" theme (ppt/theme/theme2.xml, applied to
# the slide master) is switched back to the default "Office Theme" palette.
# (ppt/theme/theme1.xml keeps the Office Theme colors that theme2.xml is
# losing - that swap happens naturally inside the host when the master's
# live theme colors are rewritten.)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$scheme = $theme.ThemeColorScheme

# ThemeColorScheme index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
# 11 hlink, 12 folHlink. RGB is the usual R + G*256 + B*65536 packed value.
$scheme.Item(3).RGB  = 6968388   # dk2      44546A
$scheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407     # accent4  FFC000
$scheme.Item(9).RGB  = 12874308  # accent5  4472C4
$scheme.Item(10).RGB = 4697456   # accent6  70AD47
$scheme.Item(11).RGB = 12673797  # hlink    0563C1
$scheme.Item(12).RGB = 7491477   # folHlink 954F72
